$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.58
$ws.Range("AA2").Value = 120
$ws.Range("AE2").Value = 50
$ws.Range("AF2").Value = 32
$ws.Range("AI2").Value = 60
$ws.Range("AJ2").Value = 110
$ws.Range("AK2").Value = 44
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 200
$ws.Range("F3").Value = 1.98
$ws.Range("G3").Value = 2.16
$ws.Range("H3").Value = 3.8
$ws.Range("J3").Value = 3.45
$ws.Range("P3").Value = 1.9
$ws.Range("Q3").Value = 1.91
$ws.Range("R3").Value = 1.35
$ws.Range("S3").Value = 3.35
$ws.Range("T3").Value = 1.75
$ws.Range("U3").Value = 2.08
$ws.Range("V3").Value = 1.29
$ws.Range("W3").Value = 1.86
$ws.Range("AA3").Value = 440
$ws.Range("AG3").Value = 10.5
$ws.Range("AL3").Value = 38
$ws.Range("AO3").Value = 55
$ws.Range("F4").Value = 1.54
$ws.Range("G4").Value = 1.63
$ws.Range("H4").Value = 7.4
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 3.75
$ws.Range("K4").Value = 4.4
$ws.Range("L4").Value = 1.43
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 3.3
$ws.Range("O4").Value = 1.35
$ws.Range("P4").Value = 1.79
$ws.Range("R4").Value = 1.29
$ws.Range("S4").Value = 3.7
$ws.Range("T4").Value = 2.04
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.12
$ws.Range("W4").Value = 2.56
$ws.Range("X4").Value = 14.5
$ws.Range("AA4").Value = 310
$ws.Range("AB4").Value = 7.2
$ws.Range("AC4").Value = 9.4
$ws.Range("AD4").Value = 70
$ws.Range("AF4").Value = 9.4
$ws.Range("AG4").Value = 9.800000000000001
$ws.Range("AH4").Value = 970
$ws.Range("AM4").Value = 580
$ws.Range("AN4").Value = 12
$ws.Range("AO4").Value = 240
$ws.Range("F5").Value = 1.38
$ws.Range("G5").Value = 1.56
$ws.Range("H5").Value = 5.7
$ws.Range("I5").Value = 9.6
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 7.6
$ws.Range("L5").Value = 1.25
$ws.Range("N5").Value = 1.02
$ws.Range("O5").Value = 1.15
$ws.Range("P5").Value = 2.44
$ws.Range("Q5").Value = 1.4
$ws.Range("R5").Value = 1.61
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.64
$ws.Range("V5").Value = 1.13
$ws.Range("W5").Value = 2.78
$ws.Range("Z5").Value = 170
$ws.Range("AB5").Value = 500
$ws.Range("AC5").Value = 42
$ws.Range("AE5").Value = 200
$ws.Range("AG5").Value = 23
$ws.Range("AI5").Value = 190
$ws.Range("AJ5").Value = 970
$ws.Range("AK5").Value = 970
$ws.Range("AL5").Value = 970
$ws.Range("AN5").Value = 5.4
$ws.Range("AO5").Value = 65
$ws.Range("F6").Value = 2.92
$ws.Range("G6").Value = 3.2
$ws.Range("H6").Value = 2.38
$ws.Range("J6").Value = 3.5
$ws.Range("L6").Value = 1.4
$ws.Range("Q6").Value = 1.9
$ws.Range("T6").Value = 1.75
$ws.Range("W6").Value = 1.45
$ws.Range("Z6").Value = 16.5
$ws.Range("AA6").Value = 130
$ws.Range("AJ6").Value = 200
$ws.Range("AL6").Value = 110
$ws.Range("AO6").Value = 22
$ws.Range("F7").Value = 1.84
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 4.4
$ws.Range("I7").Value = 4.5
$ws.Range("Q7").Value = 1.59
$ws.Range("S7").Value = 2.48
$ws.Range("U7").Value = 2.62
$ws.Range("V7").Value = 1.28
$ws.Range("W7").Value = 2.16
$ws.Range("AN7").Value = 8
$ws.Range("I8").Value = 1.29
$ws.Range("L8").Value = 1.23
$ws.Range("P8").Value = 2.52
$ws.Range("T8").Value = 2.06
$ws.Range("U8").Value = 1.78
$ws.Range("X8").Value = 32
$ws.Range("Y8").Value = 10.5
$ws.Range("AA8").Value = 10
$ws.Range("AB8").Value = 100
$ws.Range("AC8").Value = 17
$ws.Range("AE8").Value = 13.5
$ws.Range("AF8").Value = 170
$ws.Range("AH8").Value = 36
$ws.Range("AI8").Value = 40
$ws.Range("AJ8").Value = 740
$ws.Range("AL8").Value = 210
$ws.Range("AN8").Value = 290
$ws.Range("F9").Value = 1.78
$ws.Range("G9").Value = 1.85
$ws.Range("K9").Value = 4.5
$ws.Range("L9").Value = 1.26
$ws.Range("N9").Value = 5.1
$ws.Range("P9").Value = 2.4
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 1.55
$ws.Range("T9").Value = 1.61
$ws.Range("U9").Value = 2.36
$ws.Range("AO9").Value = 1000
$ws.Range("L10").Value = 1.39
$ws.Range("Q10").Value = 1.94
$ws.Range("U10").Value = 2.1
$ws.Range("AN10").Value = 11.5
$ws.Range("F11").Value = 3.6
$ws.Range("L11").Value = 1.34
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = 1.24
$ws.Range("P11").Value = 2.3
$ws.Range("Q11").Value = 1.72
$ws.Range("R11").Value = 1.54
$ws.Range("S11").Value = 2.78
$ws.Range("V11").Value = 1.84
$ws.Range("Y11").Value = 12.5
$ws.Range("Z11").Value = 14.5
$ws.Range("AB11").Value = 17
$ws.Range("AG11").Value = 14.5
$ws.Range("AL11").Value = 40
$ws.Range("AN11").Value = 28
$ws.Range("AO11").Value = 12
$ws.Range("F13").Value = 1.39
$ws.Range("Q13").Value = 1.49
$ws.Range("U13").Value = 2.24
$ws.Range("W13").Value = 3.5
$ws.Range("AA13").Value = 250
$ws.Range("AM13").Value = 95
$ws.Range("AO13").Value = 95
$ws.Range("Z14").Value = 210
$ws.Range("F15").Value = 2.96
$ws.Range("L15").Value = 1.29
$ws.Range("AE15").Value = 25
$ws.Range("AI15").Value = 85
$ws.Range("AL15").Value = 160
$ws.Range("AO15").Value = 16.5
$ws.Range("K16").Value = 7.6
$ws.Range("Q16").Value = 1.43
$ws.Range("T16").Value = 1.81
$ws.Range("U16").Value = 1.99
$ws.Range("AL16").Value = 32
$ws.Range("AN16").Value = 3.95
$ws.Range("F17").Value = 1.81
$ws.Range("H17").Value = 4
$ws.Range("N17").Value = 3.25
$ws.Range("Q17").Value = 1.56
$ws.Range("N18").Value = 1.1
$ws.Range("P18").Value = 1.94
$ws.Range("Q18").Value = 1.64
$ws.Range("L19").Value = 1.31
$ws.Range("X19").Value = 970
$ws.Range("Y19").Value = 44
$ws.Range("AN19").Value = 55
$ws.Range("G20").Value = 2.1
$ws.Range("H20").Value = 3.95
$ws.Range("I20").Value = 4.4
$ws.Range("J20").Value = 3.5
$ws.Range("T20").Value = 1.81
$ws.Range("V20").Value = 1.3
$ws.Range("W20").Value = 1.9
$ws.Range("AA20").Value = 900
$ws.Range("AE20").Value = 150
$ws.Range("G21").Value = 610
$ws.Range("I21").Value = 870
$ws.Range("N21").Value = 1.1
$ws.Range("T21").Value = 1.04
$ws.Range("U21").Value = 1.04
$ws.Range("X21").Value = 970
$ws.Range("Y21").Value = 970
$ws.Range("AB21").Value = 970
$ws.Range("AC21").Value = 970
$ws.Range("AD21").Value = 970
$ws.Range("AG21").Value = 970
$ws.Range("AH21").Value = 970
$ws.Range("J22").Value = 3.75
$ws.Range("P22").Value = 2.28
$ws.Range("Q22").Value = 1.71
$ws.Range("F23").Value = 1.67
$ws.Range("G23").Value = 1.71
$ws.Range("H23").Value = 5.1
$ws.Range("I23").Value = 5.7
$ws.Range("P23").Value = 2.34
$ws.Range("Q23").Value = 1.65
$ws.Range("R23").Value = 1.52
$ws.Range("V23").Value = 1.21
$ws.Range("W23").Value = 2.4
$ws.Range("Z23").Value = 46
$ws.Range("AA23").Value = 150
$ws.Range("AE23").Value = 65
